# skostr_hoyde.xlsx — legg til ny datarad (skostr 43, hoyde 181) i rad 7
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 43
$ws.Range("B7").Value = 181

# Flytt den aktive markeringen slik den står etter ny rad er lagt inn
$ws.Range("C7").Select()
